$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.081.51'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.323.62'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.90%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.72%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.66%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("E8").Value = '  -4.19%  '

$ws.Range("E9").Value = '  -2.25%  '

$ws.Range("E10").Value = '  -2.13%  '

$ws.Range("E11").Value = '  +0.84%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.338'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.89'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.732.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.07%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.047.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.50%  '

$ws.Range("E16").Value = '  -2.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.316.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.94%  '

$ws.Range("E18").Value = '  -3.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '318.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.84%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.78%  '

$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '63.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.170'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.28%  '

$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.04%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0740'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.89'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.70%  '

$ws.Range("E32").Value = '  +4.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.393'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.14%  '

$ws.Range("E34").Value = '  +0.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.88'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.79%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.13%  '

$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.04'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.65%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.56'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.03%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '38.56'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '307.69'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '143.84'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0952'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0500'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.560'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.59%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0214'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.12%  '

$ws.Range("E50").Value = '  +0.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.936'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.86%  '
